$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$hf = $m.HeadersFooters
$dt = $hf.DateAndTime
"Before Text=[" + $dt.Text + "] UseFormat=" + $dt.UseFormat + " Format=" + $dt.Format + " Visible=" + $dt.Visible
$dt.Text = "6/5/25"
"After Text=[" + $dt.Text + "]"
